$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Apply all cell content updates per the diff.
# D-column numeric-looking text values are forced to stay text
# (NumberFormat "@" + Style reset) so Excel does not reinterpret
# strings like "22.60" or "30.015.35" as numbers.

$c = $ws.Cells.Item(2, 4)
$c.NumberFormat = "@"
$c.Value = "30.015.35"
$c.Style = "Normal"
$ws.Cells.Item(2, 5).Value = "  -0.21%  "

$c = $ws.Cells.Item(3, 4)
$c.NumberFormat = "@"
$c.Value = "1.868.09"
$c.Style = "Normal"
$ws.Cells.Item(3, 5).Value = "  -2.87%  "

$c = $ws.Cells.Item(4, 4)
$c.NumberFormat = "@"
$c.Value = "1.002"
$c.Style = "Normal"
$ws.Cells.Item(4, 5).Value = "  +0.07%  "

$c = $ws.Cells.Item(5, 4)
$c.NumberFormat = "@"
$c.Value = "317.91"
$c.Style = "Normal"
$ws.Cells.Item(5, 5).Value = "  -2.21%  "

$ws.Cells.Item(6, 5).Value = "  +0.02%  "

$ws.Cells.Item(7, 5).Value = "  -1.73%  "

$c = $ws.Cells.Item(8, 4)
$c.NumberFormat = "@"
$c.Value = "0.3903"
$c.Style = "Normal"
$ws.Cells.Item(8, 5).Value = "  -2.42%  "

$c = $ws.Cells.Item(9, 4)
$c.NumberFormat = "@"
$c.Value = "0.08149"
$c.Style = "Normal"
$ws.Cells.Item(9, 5).Value = "  -3.69%  "

$ws.Cells.Item(10, 5).Value = "  -2.13%  "

$c = $ws.Cells.Item(11, 4)
$c.NumberFormat = "@"
$c.Value = "1.087"
$c.Style = "Normal"
$ws.Cells.Item(11, 5).Value = "  -3.13%  "

$c = $ws.Cells.Item(12, 4)
$c.NumberFormat = "@"
$c.Value = "22.60"
$c.Style = "Normal"
$ws.Cells.Item(12, 5).Value = "  +6.76%  "

$c = $ws.Cells.Item(13, 4)
$c.NumberFormat = "@"
$c.Value = "1.863.62"
$c.Style = "Normal"
$ws.Cells.Item(13, 5).Value = "  -2.92%  "

$c = $ws.Cells.Item(14, 4)
$c.NumberFormat = "@"
$c.Value = "6.248"
$c.Style = "Normal"
$ws.Cells.Item(14, 5).Value = "  -1.10%  "

$ws.Cells.Item(15, 5).Value = "  -2.81%  "

$c = $ws.Cells.Item(16, 4)
$c.NumberFormat = "@"
$c.Value = "1.002"
$c.Style = "Normal"
$ws.Cells.Item(16, 5).Value = "  +0.12%  "

$c = $ws.Cells.Item(17, 4)
$c.NumberFormat = "@"
$c.Value = "91.50"
$c.Style = "Normal"
$ws.Cells.Item(17, 5).Value = "  -3.06%  "

$ws.Cells.Item(18, 5).Value = "  -3.82%  "

$c = $ws.Cells.Item(19, 4)
$c.NumberFormat = "@"
$c.Value = "0.06326"
$c.Style = "Normal"
$ws.Cells.Item(19, 5).Value = "  -6.39%  "

$c = $ws.Cells.Item(20, 4)
$c.NumberFormat = "@"
$c.Value = "17.78"
$c.Style = "Normal"
$ws.Cells.Item(20, 5).Value = "  -1.12%  "

$ws.Cells.Item(21, 5).Value = "  -0.04%  "

$c = $ws.Cells.Item(22, 4)
$c.NumberFormat = "@"
$c.Value = "29.994.58"
$c.Style = "Normal"
$ws.Cells.Item(22, 5).Value = "  -0.28%  "

$ws.Cells.Item(23, 5).Value = "  -4.45%  "

$ws.Cells.Item(24, 5).Value = "  -1.43%  "

$c = $ws.Cells.Item(25, 4)
$c.NumberFormat = "@"
$c.Value = "2.203"
$c.Style = "Normal"
$ws.Cells.Item(25, 5).Value = "  +0.07%  "

$c = $ws.Cells.Item(26, 4)
$c.NumberFormat = "@"
$c.Value = "2.086.76"
$c.Style = "Normal"
$ws.Cells.Item(26, 5).Value = "  -2.48%  "

$c = $ws.Cells.Item(27, 4)
$c.NumberFormat = "@"
$c.Value = "160.53"
$c.Style = "Normal"
$ws.Cells.Item(27, 5).Value = "  +0.25%  "

$c = $ws.Cells.Item(28, 4)
$c.NumberFormat = "@"
$c.Value = "20.79"
$c.Style = "Normal"
$ws.Cells.Item(28, 5).Value = "  -1.14%  "

$c = $ws.Cells.Item(29, 4)
$c.NumberFormat = "@"
$c.Value = "2.212"
$c.Style = "Normal"
$ws.Cells.Item(29, 5).Value = "  -10.01%  "

$ws.Cells.Item(30, 5).Value = "  -2.43%  "

$c = $ws.Cells.Item(31, 4)
$c.NumberFormat = "@"
$c.Value = "0.1028"
$c.Style = "Normal"
$ws.Cells.Item(31, 5).Value = "  -2.64%  "

$c = $ws.Cells.Item(32, 4)
$c.NumberFormat = "@"
$c.Value = "1.035"
$c.Style = "Normal"
$ws.Cells.Item(32, 5).Value = "  -3.76%  "

$c = $ws.Cells.Item(33, 4)
$c.NumberFormat = "@"
$c.Value = "5.850"
$c.Style = "Normal"
$ws.Cells.Item(33, 5).Value = "  -3.64%  "

$c = $ws.Cells.Item(34, 4)
$c.NumberFormat = "@"
$c.Value = "3.725"
$c.Style = "Normal"
$ws.Cells.Item(34, 5).Value = "  +1.83%  "

$ws.Cells.Item(35, 5).Value = "  -3.43%  "

$c = $ws.Cells.Item(36, 4)
$c.NumberFormat = "@"
$c.Value = "5.162"
$c.Style = "Normal"
$ws.Cells.Item(36, 5).Value = "  -0.57%  "

$c = $ws.Cells.Item(37, 4)
$c.NumberFormat = "@"
$c.Value = "0.06308"
$c.Style = "Normal"
$ws.Cells.Item(37, 5).Value = "  -4.43%  "

$c = $ws.Cells.Item(38, 4)
$c.NumberFormat = "@"
$c.Value = "0.2130"
$c.Style = "Normal"
$ws.Cells.Item(38, 5).Value = "  -3.68%  "

$ws.Cells.Item(39, 5).Value = "  -6.10%  "

$c = $ws.Cells.Item(40, 4)
$c.NumberFormat = "@"
$c.Value = "8.483"
$c.Style = "Normal"
$ws.Cells.Item(40, 5).Value = "  -5.78%  "

$ws.Cells.Item(41, 2).Value = "TheSandbox"
$ws.Cells.Item(41, 3).Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$c = $ws.Cells.Item(41, 4)
$c.NumberFormat = "@"
$c.Value = "0.6242"
$c.Style = "Normal"
$ws.Cells.Item(41, 5).Value = "  -4.31%  "

$ws.Cells.Item(42, 2).Value = "TrustWalletToken"
$ws.Cells.Item(42, 3).Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$c = $ws.Cells.Item(42, 4)
$c.NumberFormat = "@"
$c.Value = "1.205"
$c.Style = "Normal"
$ws.Cells.Item(42, 5).Value = "  -2.78%  "

$c = $ws.Cells.Item(43, 4)
$c.NumberFormat = "@"
$c.Value = "11.17"
$c.Style = "Normal"
$ws.Cells.Item(43, 5).Value = "  -2.18%  "

$ws.Cells.Item(44, 5).Value = "  -0.10%  "

$c = $ws.Cells.Item(45, 4)
$c.NumberFormat = "@"
$c.Value = "0.5847"
$c.Style = "Normal"
$ws.Cells.Item(45, 5).Value = "  -4.60%  "

$c = $ws.Cells.Item(46, 4)
$c.NumberFormat = "@"
$c.Value = "12.74"
$c.Style = "Normal"
$ws.Cells.Item(46, 5).Value = "  -2.92%  "

$c = $ws.Cells.Item(47, 4)
$c.NumberFormat = "@"
$c.Value = "3.624"
$c.Style = "Normal"
$ws.Cells.Item(47, 5).Value = "  -3.25%  "

$c = $ws.Cells.Item(48, 4)
$c.NumberFormat = "@"
$c.Value = "1.980"
$c.Style = "Normal"
$ws.Cells.Item(48, 5).Value = "  -3.64%  "

$c = $ws.Cells.Item(49, 4)
$c.NumberFormat = "@"
$c.Value = "121.82"
$c.Style = "Normal"
$ws.Cells.Item(49, 5).Value = "  -2.70%  "

$c = $ws.Cells.Item(50, 4)
$c.NumberFormat = "@"
$c.Value = "1.199"
$c.Style = "Normal"
$ws.Cells.Item(50, 5).Value = "  -3.35%  "

$c = $ws.Cells.Item(51, 4)
$c.NumberFormat = "@"
$c.Value = "1.136"
$c.Style = "Normal"
$ws.Cells.Item(51, 5).Value = "  -1.08%  "
